$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 134. This shifts the existing rows
# 134-151 down to 135-152 (and the sheet dimension grows to A1:T152).
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new data record.
$ws.Range("A134").Value = 6
$ws.Range("B134").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C134").Value = "Metropolitana"
$ws.Range("D134").Value = 44449
$ws.Range("E134").Value = 13
$ws.Range("F134").Value = "Fruta"
$ws.Range("G134").Value = 100101
$ws.Range("H134").Value = "Berries"
$ws.Range("I134").Value = 100101001
$ws.Range("J134").Value = "Arándano (blue)"
$ws.Range("K134").Value = "Sin especificar"
$ws.Range("L134").Value = "Segunda"
$ws.Range("M134").Value = 1350
$ws.Range("N134").Value = 10000
$ws.Range("O134").Value = 10000
$ws.Range("P134").Value = 10000
$ws.Range("Q134").Value = "$/bandeja 2 kilos"
$ws.Range("R134").Value = "Provincia de Linares"
$ws.Range("S134").Value = 5000
$ws.Range("T134").Value = 2
